$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.808.54"
$ws.Range("E2").Value = "  +1.58%  "
$ws.Range("D3").Value = "3.485.62"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.79"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.33"
$ws.Range("E6").Value = "  +2.28%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.484.46"
$ws.Range("E8").Value = "  +0.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.588"
$ws.Range("E9").Value = "  +4.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.29"
$ws.Range("E10").Value = "  -3.99%  "
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.443"
$ws.Range("E12").Value = "  -1.99%  "
$ws.Range("D13").Value = "4.082.77"
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("E14").Value = "  -1.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000196"
$ws.Range("E15").Value = "  -2.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.84"
$ws.Range("E16").Value = "  +2.98%  "
$ws.Range("D17").Value = "65.758.97"
$ws.Range("E17").Value = "  +1.54%  "
$ws.Range("D18").Value = "3.432.11"
$ws.Range("E18").Value = "  -0.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.43"
$ws.Range("E19").Value = "  -0.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.31"
$ws.Range("E20").Value = "  -0.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "391.58"
$ws.Range("E22").Value = "  -3.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.551"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.58"
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000125"
$ws.Range("E26").Value = "  +1.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.74"
$ws.Range("E28").Value = "  -0.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("E30").Value = "  +6.46%  "
$ws.Range("E31").Value = "  +3.06%  "
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.72"
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("E34").Value = "  -3.89%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  +0.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.55"
$ws.Range("E37").Value = "  +3.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "163.41"
$ws.Range("E38").Value = "  +1.58%  "
$ws.Range("E39").Value = "  +4.07%  "
$ws.Range("D40").Value = "3.079.03"
$ws.Range("E40").Value = "  +5.79%  "
$ws.Range("E41").Value = "  -2.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.31"
$ws.Range("E42").Value = "  -1.86%  "
$ws.Range("E43").Value = "  -2.13%  "
$ws.Range("E44").Value = "  +1.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "43.02"
$ws.Range("E45").Value = "  +1.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.778"
$ws.Range("E46").Value = "  +0.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.68"
$ws.Range("E47").Value = "  +7.44%  "
$ws.Range("E48").Value = "  +2.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.25"
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "311.21"
$ws.Range("E50").Value = "  +4.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.69"
$ws.Range("E51").Value = "  +1.35%  "
